$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values per diff (sheet "Sheet1", data table with header row 1)
$ws.Range("C3").Value = -12.2348
$ws.Range("A9").Value = -20.33529999999998
$ws.Range("A18").Value = -23.01420000000002
$ws.Range("A20").Value = -22.17760000000003
